# Armenia Premier League 2023-2024 sheet update
# 1) Swap the F:V content of rows 55 and 56 (the two matches played on the
#    same date got reordered upstream).
# 2) Append 8 new match rows (85-92) with the same column layout/styling as
#    the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap rows 55 and 56 (columns F through V only; A:E stay identical)
# ---------------------------------------------------------------------
$row55 = $ws.Range("F55:V55").Value()
$row56 = $ws.Range("F56:V56").Value()

$ws.Range("F55:V55").Value = $row56
$ws.Range("F56:V56").Value = $row55

# ---------------------------------------------------------------------
# 2) Append rows 85-92 with the same formatting as row 84
# ---------------------------------------------------------------------
$ws.Range("A84:V84").Copy()
$ws.Range("A85:V92").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row=85; A=84; E=45264.54166666666;  F="Noah";           G=3; H="Van";            I=1; J=1.03; K="04/12/2023 12:19"; L=1.03; M="04/12/2023 12:19"; N=13.1; O="04/12/2023 12:24"; P=13.1; Q="04/12/2023 12:24"; R=18.22; S="04/12/2023 12:24"; T=18.22; U="04/12/2023 12:24"; V="https://www.betexplorer.com/football/armenia/premier-league/noah-van/b9qqkl3N/" },
    @{ Row=86; A=85; E=45264.625;          F="BKMA";           G=2; H="Urartu";         I=0; J=4.99; K="03/12/2023 03:12"; L=5.3;  M="04/12/2023 14:37"; N=4.1;  O="03/12/2023 03:12"; P=4.39; Q="04/12/2023 14:37"; R=1.54;  S="03/12/2023 03:12"; T=1.56;  U="04/12/2023 14:37"; V="https://www.betexplorer.com/football/armenia/premier-league/bkma-urartu/KIpujUmH/" },
    @{ Row=87; A=86; E=45265.54166666666;  F="Ararat Yerevan";  G=0; H="Ararat-Armenia"; I=1; J=6.74; K="04/12/2023 01:12"; L=8.58; M="05/12/2023 12:59"; N=4.35; O="04/12/2023 01:12"; P=4.72; Q="05/12/2023 12:54"; R=1.4;   S="04/12/2023 01:12"; T=1.37;  U="05/12/2023 12:59"; V="https://www.betexplorer.com/football/armenia/premier-league/ararat-yerevan-ararat-armenia/WdEc8B3o/" },
    @{ Row=88; A=87; E=45266.45833333334;  F="Alashkert";       G=1; H="Pyunik Yerevan"; I=1; J=4.25; K="04/12/2023 23:12"; L=3.47; M="06/12/2023 10:58"; N=3.58; O="04/12/2023 23:12"; P=3.39; Q="06/12/2023 10:58"; R=1.73;  S="04/12/2023 23:12"; T=2.11;  U="06/12/2023 10:58"; V="https://www.betexplorer.com/football/armenia/premier-league/alashkert-pyunik-yerevan/UizFqSQp/" },
    @{ Row=89; A=88; E=45268.625;          F="BKMA";           G=0; H="Noah";           I=2; J=8.19; K="07/12/2023 03:12"; L=7.5;  M="08/12/2023 14:53"; N=5.02; O="07/12/2023 03:12"; P=4.7;  Q="08/12/2023 14:53"; R=1.3;   S="07/12/2023 03:12"; T=1.41;  U="08/12/2023 14:53"; V="https://www.betexplorer.com/football/armenia/premier-league/bkma-noah/CC9I3mIG/" },
    @{ Row=90; A=89; E=45269.45833333334;  F="Van";             G=4; H="Shirak Gyumri"; I=0; J=3.22; K="07/12/2023 23:12"; L=3.66; M="09/12/2023 10:59"; N=3.2;  O="07/12/2023 23:12"; P=3.54; Q="09/12/2023 10:58"; R=2.14;  S="07/12/2023 23:12"; T=1.82;  U="09/12/2023 10:59"; V="https://www.betexplorer.com/football/armenia/premier-league/van-shirak-gyumri/0x6E4T2A/" },
    @{ Row=91; A=90; E=45269.625;          F="Urartu";          G=1; H="Ararat-Armenia"; I=3; J=3.22; K="08/12/2023 03:13"; L=3.41; M="09/12/2023 14:52"; N=3.35; O="08/12/2023 03:13"; P=3.55; Q="09/12/2023 14:52"; R=2.07;  S="08/12/2023 03:13"; T=2.07;  U="09/12/2023 14:52"; V="https://www.betexplorer.com/football/armenia/premier-league/urartu-ararat-armenia/AmD17VIi/" },
    @{ Row=92; A=91; E=45271.45833333334;  F="Pyunik Yerevan";  G=2; H="Ararat Yerevan"; I=2; J=1.22; K="09/12/2023 23:12"; L=1.2;  M="11/12/2023 10:56"; N=5.89; O="09/12/2023 23:12"; P=7.08; Q="11/12/2023 10:58"; R=9.46;  S="09/12/2023 23:12"; T=12.25; U="11/12/2023 10:58"; V="https://www.betexplorer.com/football/armenia/premier-league/pyunik-yerevan-ararat-yerevan/vNH56kYc/" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value  = $r.A
    $ws.Cells.Item($rowNum, 2).Value  = "armenia"
    $ws.Cells.Item($rowNum, 3).Value  = "premier-league"
    $ws.Cells.Item($rowNum, 4).Value  = "2023-2024"
    $ws.Cells.Item($rowNum, 5).Value  = $r.E
    $ws.Cells.Item($rowNum, 6).Value  = $r.F
    $ws.Cells.Item($rowNum, 7).Value  = $r.G
    $ws.Cells.Item($rowNum, 8).Value  = $r.H
    $ws.Cells.Item($rowNum, 9).Value  = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
    $ws.Cells.Item($rowNum, 21).Value = $r.U
    $ws.Cells.Item($rowNum, 22).Value = $r.V
}
